$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 111: new timelog entry (2015-05-01, 23:00 -> 24:00) ---
$ws.Cells.Item(111, 1).Value = 42125
$ws.Cells.Item(111, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(111, 2).Value = 0.95833333333333337
$ws.Cells.Item(111, 2).NumberFormat = "h:mm"
$ws.Cells.Item(111, 3).Value = 1
$ws.Cells.Item(111, 3).NumberFormat = "h:mm"
$ws.Cells.Item(111, 5).Value = "thesis chapter state of the art"

# --- Row 112: new timelog entry (2015-05-02, 00:00 -> 02:30) ---
$ws.Cells.Item(112, 1).Value = 42126
$ws.Cells.Item(112, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(112, 2).Value = 0
$ws.Cells.Item(112, 2).NumberFormat = "h:mm"
$ws.Cells.Item(112, 3).Value = 0.10416666666666667
$ws.Cells.Item(112, 3).NumberFormat = "h:mm"
$ws.Cells.Item(112, 5).Value = "thesis chapter state of the art"

# --- Rows 113-124: give the still-empty A/B/C cells the same formatting
# as the rest of the table (they previously had no cells at all there) ---
for ($r = 113; $r -le 124; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $ws.Cells.Item($r, 2).NumberFormat = "h:mm"
    $ws.Cells.Item($r, 3).NumberFormat = "h:mm"
}

# --- View state: scroll down a little and move the selection, matching
# the author's on-screen position after adding the two rows above ---
$ws.Range("E117").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
